$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.958.82'
$ws.Range("E2").Value = '  -1.45%  '
$ws.Range("D3").Value = '3.424.49'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.42%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.432.49'
$ws.Range("E8").Value = '  -0.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.552'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.49%  '
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.119'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.423'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.13%  '
$ws.Range("D13").Value = '4.030.53'
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.135'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000173'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.65%  '
$ws.Range("D17").Value = '64.020.19'
$ws.Range("E17").Value = '  -1.32%  '
$ws.Range("D18").Value = '3.446.22'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '377.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.76'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.23'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.516'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000116'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.53'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.41%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.40'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.93'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.14'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("E37").Value = '  +11.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.81'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.59%  '
$ws.Range("D39").Value = '2.801.62'
$ws.Range("E39").Value = '  -3.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0726'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.80'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.90'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.48%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.89%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '26.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.40'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0306'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '329.20'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.34%  '
$ws.Range("E51").Value = '  -2.79%  '
